$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8 and 9 entirely (shifts old rows 10-14 up to become rows 8-12)
$ws.Rows("8:9").Delete()

# Fix up the values in the (new) rows 8 and 9 to match the target state
$ws.Range("A8").Value = "TestLoad"
$ws.Range("B8").Value = "Template"
$ws.Range("C8").Value = "No"
$ws.Range("D8").Value = "Test"
$ws.Range("E8").Value = "Souce5"
$ws.Range("F8").Value = "Souce5"
$ws.Range("G8").Value = "Souce5"
$ws.Range("H8").Value = "Souce5"
$ws.Range("I8").Value = "bengteth\administrator"
$ws.Range("J8").Value = "123"

$ws.Range("A9").Value = "24/10/2019 14:00:41"
$ws.Range("B9").Value = "TestSM"
$ws.Range("C9").Value = "Yes"
$ws.Range("D9").Value = "SMS Testing"
$ws.Range("E9").Value = "18"
$ws.Range("F9").Value = "18"
$ws.Range("G9").Value = "ReTesting"
$ws.Range("H9").Value = "18"
$ws.Range("I9").Value = "bengteth\administrator"
$ws.Range("J9").Value = "24/12/2019 15:24:53"

Write-Host "Done"
